$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '52.346.99'
$ws.Range("E2").Value = '  +1.45%  '
$ws.Range("D3").Value = '2.907.57'
$ws.Range("E3").Value = '  +4.17%  '
$ws.Range("E4").Value = '  +0.13%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '353.73'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = '  +0.22%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '114.18'
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = '  +2.73%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.559'
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").Value = '  +1.05%  '
$ws.Range("E8").Value = '  +0.06%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.626'
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = '  -0.32%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '40.19'
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = '  +0.47%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0864'
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = '  +3.33%  '
$ws.Range("E12").Value = '  +0.62%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '19.86'
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E14").Value = '  +0.63%  '
$ws.Range("D15").Value = '3.367.10'
$ws.Range("E15").Value = '  +4.26%  '
$ws.Range("E16").Value = '  +6.16%  '
$ws.Range("D17").Value = '2.919.50'
$ws.Range("E17").Value = '  +4.28%  '
$ws.Range("D18").Value = '52.383.11'
$ws.Range("E18").Value = '  +1.61%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.66'
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = '  +1.13%  '
$ws.Range("E20").Value = '  +3.55%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.16'
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = '  +4.71%  '
$ws.Range("E22").Value = '  +1.01%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '70.89'
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = '  +0.97%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '269.39'
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = '  +0.93%  '
$ws.Range("E25").Value = '  +1.76%  '
$ws.Range("E26").Value = '  +8.25%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '26.85'
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = '  +3.10%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.00'
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = '  -0.06%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '10.64'
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = '  +3.09%  '
$ws.Range("B30").Value = 'Hedera'
$ws.Range("C30").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.102'
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = '  +13.88%  '
$ws.Range("B31").Value = 'InjectiveProtocol'
$ws.Range("C31").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '37.83'
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = '  -2.76%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.58'
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = '  +7.79%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.27'
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = '  +12.59%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '53.37'
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = '  +1.56%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0451'
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = '  -0.34%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.97'
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = '  -12.74%  '
$ws.Range("E37").Value = '  +0.00%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.35'
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = '  +6.30%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '18.97'
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = '  +0.91%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.06'
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = '  +3.16%  '
$ws.Range("E41").Value = '  +11.51%  '
$ws.Range("E42").Value = '  +2.16%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '23.17'
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = '  +6.67%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '120.56'
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = '  -0.56%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.61'
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = '  +6.18%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.20'
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = '  -2.16%  '
$ws.Range("E47").Value = '  +4.00%  '
$ws.Range("D48").Value = '2.184.18'
$ws.Range("E48").Value = '  +3.81%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.263'
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = '  +20.97%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0351'
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = '  +14.41%  '
$ws.Range("E51").Value = '  +0.45%  '
